$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.765038666666667
$ws.Range("H2").Value = 8.295116
$ws.Range("I2").Value = 0.5643238178805899
$ws.Range("J2").Value = 0.5643238178805899
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1375686666666667
$ws.Range("N2").Value = 0.412706
$ws.Range("O2").Value = 0.2896572731203081
$ws.Range("P2").Value = 0.2896572731203081
$ws.Range("Q2").Value = 0.3803826826551112
$ws.Range("R2").Value = 3.423444143896
$ws.Range("S2").Value = 0.1634604982441331
$ws.Range("T2").Value = 0.1634604982441331

$ws.Range("G3").Value = 2.765038666666667
$ws.Range("H3").Value = 8.295116
$ws.Range("I3").Value = 0.5643238178805899
$ws.Range("J3").Value = 0.5643238178805899
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.3373673333333334
$ws.Range("N3").Value = 1.012102
$ws.Range("O3").Value = 0.7103427268796919
$ws.Range("P3").Value = 0.7103427268796919
$ws.Range("Q3").Value = 0.932833721536889
$ws.Range("R3").Value = 8.395503493832001
$ws.Range("S3").Value = 0.4008633196364569
$ws.Range("T3").Value = 0.4008633196364569

$ws.Range("G4").Value = 2.134699
$ws.Range("H4").Value = 6.404097
$ws.Range("I4").Value = 0.4356761821194101
$ws.Range("J4").Value = 0.4356761821194101
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.1375686666666667
$ws.Range("N4").Value = 0.412706
$ws.Range("O4").Value = 0.2896572731203081
$ws.Range("P4").Value = 0.2896572731203081
$ws.Range("Q4").Value = 0.2936676951646667
$ws.Range("R4").Value = 2.643009256482
$ws.Range("S4").Value = 0.1261967748761751
$ws.Range("T4").Value = 0.1261967748761751

$ws.Range("G5").Value = 2.134699
$ws.Range("H5").Value = 6.404097
$ws.Range("I5").Value = 0.4356761821194101
$ws.Range("J5").Value = 0.4356761821194101
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.3373673333333334
$ws.Range("N5").Value = 1.012102
$ws.Range("O5").Value = 0.7103427268796919
$ws.Range("P5").Value = 0.7103427268796919
$ws.Range("Q5").Value = 0.7201777090993333
$ws.Range("R5").Value = 6.481599381894
$ws.Range("S5").Value = 0.309479407243235
$ws.Range("T5").Value = 0.309479407243235
